$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Common")

# The engine does not re-anchor cell comments when rows are inserted (unlike
# merged cells / data validations / formulas, which do shift). So capture
# every existing comment below the insertion point (row 70 and down) up
# front, remove them, perform the row insert, and then re-create each one
# shifted down by one row.
$commentData = New-Object System.Collections.ArrayList
for ($r = 70; $r -le 129; $r++) {
    $cell = $ws.Range("A$r")
    $cm = $cell.Comment
    if ($cm -ne $null) {
        [void]$commentData.Add(@{ Row = $r; Text = $cm.Text() })
        $cm.Delete()
    }
}

# Insert the new row for "VSTAT License File" above the old row 70
# ("SD-WAN Portal License File"), pushing everything down by one.
$ws.Rows.Item(70).Insert()

# The inserted row's B cell doesn't pick up the same cell style used by all
# its sibling rows (s="9"); copy formats down from the row below so B70
# matches the rest of the "value" column.
$ws.Range("B71").Copy() | Out-Null
$ws.Range("B70").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Populate the newly inserted row.
$ws.Range("A70").Value = "VSTAT License File"
$ws.Range("A70").AddComment("Optional License File for Elasticsearch [default: ]") | Out-Null

# Re-create the previously captured comments one row further down.
foreach ($entry in $commentData) {
    $newRow = [int]$entry.Row + 1
    $target = $ws.Range("A$newRow")
    $target.AddComment($entry.Text) | Out-Null
}
